$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 785 (1-based), shifting existing rows 785+ down by one.
$ws.Rows.Item(785).Insert()

# Populate the newly inserted row 785 with the new data point.
$ws.Cells.Item(785, 1).NumberFormat = "@"
$ws.Cells.Item(785, 1).Value = "2026/02/12"
$ws.Cells.Item(785, 2).Value = "木"
$ws.Cells.Item(785, 3).Value = 17
$ws.Cells.Item(785, 4).Value = 34
